$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8308
$ws1.Range("F5").Value = 77
$ws1.Range("F6").Value = 1101
$ws1.Range("F10").Value = 216
$ws1.Range("F11").Value = 60

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8308
$ws4.Range("F5").Value = 77
$ws4.Range("F6").Value = 1101
$ws4.Range("F11").Value = 216
$ws4.Range("F12").Value = 60
